$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Replace the text of the trailing italic paragraph (near the end
#    of the document) with the new "Prompt:" text. We rebuild the
#    Range from the Document (rather than reusing Paragraph.Range)
#    so the assignment properly overwrites the old text instead of
#    merely inserting before it.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$targetPara = $null
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd("`r")
    if ($text -eq "Discover the flashy gameplay, bonus features, and luxurious visual design of Crystal Land. Play for free or place real bets in online casinos.") {
        $targetPara = $p
        break
    }
}
$s = $targetPara.Range.Start
$e = $targetPara.Range.End
$textRange = $d.Range($s, $e)
$textRange.Text = "Prompt: Create a cartoon-style feature image for Crystal Land with a happy Maya warrior wearing glasses. The Maya warrior should be standing in front of a background of vibrant, sparkling crystals and gems, looking excited and enthusiastic about the game. He should be holding a smartphone or tablet, with the Crystal Land game displayed on the screen, and there should be a speech bubble next to him with the text `"Join the Crystal Land adventure!`" written in it. The cartoon-style image should be colorful and eye-catching, with the Maya warrior wearing modern-style glasses to give the image a modern touch."

# ------------------------------------------------------------------
# 2. Remove the trailing duplicate "Play Crystal Land..." paragraph
#    that used to sit just before the paragraph we edited above.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd("`r") -eq "Play Crystal Land for Free - A Luxurious Jewel-Themed Slot") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3. Insert a new "Meta description" paragraph right after the
#    first (Heading1) paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$ms = $metaPara.Range.Start
$me = $metaPara.Range.End

$metaRange = $d.Range($ms, $me)
$metaRange.InsertBefore("Meta description: Discover the flashy gameplay, bonus features, and luxurious visual design of Crystal Land. Play for free or place real bets in online casinos.")

# Make "Meta description" bold; leave the rest (incl. the leading
# colon) regular.
$boldEnd = $ms + ("Meta description").Length
$boldRange = $d.Range($ms, $boldEnd)
$boldRange.Bold = 1

# Leave a leading empty run in place, mirroring the empty leading
# run convention used by every other body paragraph in this document.
$leadRange = $d.Range($ms, $ms)
$leadRange.InsertBefore("")
